$wb = $excel.ActiveWorkbook

# --- Sheet: Metadata ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# URL: ncpi-fhir -> nih-ncpi
$wsMeta.Range("B2").Value = "https://nih-ncpi.github.io/ncpi-fhir-ig/ValueSet/condition-inheritance-vs"

# Version: 0.1.0 -> 0.2.0
$wsMeta.Range("B3").Value = "0.2.0"

# Date: 2022-08-23T15:00:44+00:00 -> 2022-09-13T16:54:38+00:00
$wsMeta.Range("B8").Value = "2022-09-13T16:54:38+00:00"

# --- Sheet: Include from Condition Inheri (System URI) ---
$wsInclude = $wb.Worksheets.Item("Include from Condition Inheri")
$wsInclude.Range("B4").Value = "https://nih-ncpi.github.io/ncpi-fhir-ig/CodeSystem/ConditionInheritanceMode"
